$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new key/value pairs to the language table (rows 17 and 18)
$ws.Range("A17").Value = "search_title"
$ws.Range("B17").Value = "Search: {0}"

$ws.Range("A18").Value = "search_button"
$ws.Range("B18").Value = "SEARCH"

# Update selection to match the new active cell after the edit
$ws.Range("A18").Select()
